$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) column, forcing text format to avoid numeric auto-conversion
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.270.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.112.94"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.105.41"
$ws.Range("D8").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.78"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000248"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.626.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.169.71"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.109.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "467.57"
$ws.Range("D20").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.58"
$ws.Range("D27").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.86"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0870"
$ws.Range("D34").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.36"
$ws.Range("D37").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "438.18"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.70"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.914.33"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.280"
$ws.Range("D44").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "35.10"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.999"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.23"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.60"
$ws.Range("D51").Style = "Normal"

# Update Volume(1h) (E) column
$ws.Range("E2").Value = "  +6.41%  "
$ws.Range("E3").Value = "  +4.23%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  +3.21%  "
$ws.Range("E6").Value = "  +4.88%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +4.22%  "
$ws.Range("E9").Value = "  +2.28%  "
$ws.Range("E10").Value = "  +13.96%  "
$ws.Range("E11").Value = "  +7.64%  "
$ws.Range("E12").Value = "  +3.85%  "
$ws.Range("E13").Value = "  +8.13%  "
$ws.Range("E14").Value = "  +5.09%  "
$ws.Range("E15").Value = "  +0.77%  "
$ws.Range("E16").Value = "  +4.21%  "
$ws.Range("E17").Value = "  +1.26%  "
$ws.Range("E18").Value = "  +6.19%  "
$ws.Range("E19").Value = "  +4.14%  "
$ws.Range("E20").Value = "  +6.56%  "
$ws.Range("E21").Value = "  +3.90%  "
$ws.Range("E22").Value = "  +0.53%  "
$ws.Range("E23").Value = "  +7.34%  "
$ws.Range("E24").Value = "  +0.23%  "
$ws.Range("E25").Value = "  +2.49%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  +10.79%  "
$ws.Range("E28").Value = "  -0.24%  "
$ws.Range("E29").Value = "  +4.63%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("E31").Value = "  +9.65%  "
$ws.Range("E32").Value = "  +4.43%  "
$ws.Range("E33").Value = "  +5.56%  "
$ws.Range("E34").Value = "  +12.05%  "
$ws.Range("E35").Value = "  +16.25%  "
$ws.Range("E36").Value = "  +7.03%  "
$ws.Range("E37").Value = "  +21.31%  "
$ws.Range("E38").Value = "  +3.14%  "
$ws.Range("E39").Value = "  +3.89%  "
$ws.Range("E40").Value = "  +9.31%  "
$ws.Range("E41").Value = "  +0.44%  "
$ws.Range("E42").Value = "  +6.56%  "
$ws.Range("E43").Value = "  +5.04%  "
$ws.Range("E44").Value = "  +11.23%  "
$ws.Range("E45").Value = "  +5.09%  "
$ws.Range("E46").Value = "  +8.06%  "
$ws.Range("E47").Value = "  +0.85%  "
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("E49").Value = "  +0.82%  "
$ws.Range("E50").Value = "  +1.40%  "
$ws.Range("E51").Value = "  +5.55%  "

# Row 47/48: swap Arweave/USDe, update Coin (B) and Link (C) columns
$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
